# Apply the RTSC 3-compartment data change:
# Row 1 takes on the values that were previously in Row 2 (A2:C2),
# and rows 2 through 11 are cleared to 0 (rows 12-16 were already 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the old Row 2 values before overwriting anything.
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2

# Row 1 becomes what Row 2 used to be.
$ws.Range("A1").Value = $a2
$ws.Range("B1").Value = $b2
$ws.Range("C1").Value = $c2

# Rows 2 through 11 become 0.
$ws.Range("A2:C11").Value = 0
